$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previously last row (49) becomes a normal data row (datetime number format)
$ws.Range("A49").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Add new row 50 as the new last row
$ws.Range("A50").Value = 45790
$ws.Range("B50").Value = 206
$ws.Range("C50").Value = 214
$ws.Range("D50").Value = 211

# New last row gets the date-only number format
$ws.Range("A50").NumberFormat = "YYYY-MM-DD"
